$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Field-name values copied from the highlighted cells in Sheet1 column B
# (B2, B10, B11, B13, B26, B27, B41, B54, B57, B63, B64, B76, B95-B100,
#  B129-B133, B160) into the new Sheet2, column A, rows 1-24.
$values = @(
    "study_id",
    "age",
    "race_and_ethnicity",
    "sex",
    "education",
    "income",
    "essential",
    "total_interaction",
    "social_interaction",
    "more_time_household",
    "after_covid_interaction",
    "hobbies___1",
    "physical_overall",
    "mental_overall",
    "socio_emotional_overall",
    "physical_activities",
    "covid_exercise",
    "mental_health",
    "phy_health_matrix",
    "men_health_matrix",
    "time_spent_with_family",
    "time_spent_working",
    "time_spent_on_hobbies",
    "contributed"
)

# Add the new worksheet after Sheet1 and rename it to Sheet2
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Restore the highlighted-cell selection on Sheet1 (the field names that
# were copied into Sheet2); B160 was the last / active cell of that
# multi-area selection.
[void]$ws1.Range("B2,B10,B11,B13,B26,B27,B41,B54,B57,B63,B64,B76,B95,B96,B97,B98,B99,B100,B129,B130,B131,B132,B133,B160").Select()
[void]$ws1.Range("B160").Select()

# Select the newly populated range on Sheet2 and make it the active sheet/tab
[void]$ws2.Range("A1:A24").Select()
[void]$ws2.Activate()
